$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(11, "Katarzyna", "Aubert",  "katarzyna.aubert@example.com", "female", "Switzerland", "2025-03-01 14:45:47"),
    @(12, "Jairo",     "Vergara", "jairo.vergara@example.com",    "male",   "Mexico",      "2025-03-01 14:45:47"),
    @(13, "Ananya",    "Pujari",  "ananya.pujari@example.com",    "female", "India",       "2025-03-01 14:45:47"),
    @(14, "Juanita",   "Ramos",   "juanita.ramos@example.com",    "female", "Australia",   "2025-03-01 14:45:47"),
    @(15, "Judy",      "Curtis",  "judy.curtis@example.com",      "female", "Ireland",     "2025-03-01 14:45:47")
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
